$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 452.16666
$ws.Range("I28").Value = 356.9091
$ws.Range("K28").Value = 356.9091
$ws.Range("M28").Value = 128.0909
$ws.Range("H101").Value = 358.33334
$ws.Range("I101").Value = 303.5
$ws.Range("J101").Value = 468
$ws.Range("K101").Value = 910.5
$ws.Range("L101").Value = 1404
$ws.Range("M101").Value = 711.5
$ws.Range("N101").Value = -4648
$ws.Range("H112").Value = 1875.5714
$ws.Range("J112").Value = 2546
$ws.Range("L112").Value = 7638
$ws.Range("N112").Value = -9854

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 19995
$ws.Range("J9").Value = 19995
$ws.Range("L9").Value = 19995
$ws.Range("N9").Value = -20335
$ws.Range("H20").Value = 19995
$ws.Range("J20").Value = 19995
$ws.Range("L20").Value = 19995
$ws.Range("N20").Value = -20535
$ws.Range("H32").Value = 15877162
$ws.Range("I32").Value = 18183550
$ws.Range("K32").Value = 18183550
$ws.Range("M32").Value = -18183263
$ws.Range("H37").Value = 11945
$ws.Range("I37").Value = 3890
$ws.Range("K37").Value = 3890
$ws.Range("M37").Value = -3617
$ws.Range("H44").Value = 14499
$ws.Range("I44").Value = 9998
$ws.Range("J44").Value = 19000
$ws.Range("K44").Value = 9998
$ws.Range("L44").Value = 19000
$ws.Range("M44").Value = -9510
$ws.Range("N44").Value = -19976
$ws.Range("H55").Value = 17500
$ws.Range("J55").Value = 20000
$ws.Range("L55").Value = 20000
$ws.Range("N55").Value = -20630
$ws.Range("H61").Value = 2993.1538
$ws.Range("I61").Value = 1267.8889
$ws.Range("K61").Value = 1267.8889
$ws.Range("M61").Value = -1055.8889
$ws.Range("I63").Value = 2555
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 2555
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1869
$ws.Range("N63").Value = -3872
$ws.Range("I66").Value = 2555
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 12775
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -9343
$ws.Range("N66").Value = -19364
$ws.Range("H80").Value = 25201.428
$ws.Range("J80").Value = 25201.428
$ws.Range("L80").Value = 25201.428
$ws.Range("N80").Value = -27197.428
$ws.Range("H83").Value = 25201.428
$ws.Range("J83").Value = 25201.428
$ws.Range("L83").Value = 75604.284
$ws.Range("N83").Value = -85588.284
$ws.Range("H136").Value = 2993.1538
$ws.Range("I136").Value = 1267.8889
$ws.Range("K136").Value = 3803.6667
$ws.Range("M136").Value = -1253.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 12113
$ws.Range("I82").Value = 4604.2856
$ws.Range("J82").Value = 29633.334
$ws.Range("K82").Value = 4604.2856
$ws.Range("L82").Value = 29633.334
$ws.Range("M82").Value = -4221.2856
$ws.Range("N82").Value = -30399.334
$ws.Range("H85").Value = 12113
$ws.Range("I85").Value = 4604.2856
$ws.Range("J85").Value = 29633.334
$ws.Range("K85").Value = 4604.2856
$ws.Range("L85").Value = 29633.334
$ws.Range("M85").Value = -3278.2856
$ws.Range("N85").Value = -32285.334
$ws.Range("H134").Value = 2923.6785
$ws.Range("I134").Value = 2588.182
$ws.Range("J134").Value = 4153.8335
$ws.Range("K134").Value = 7764.545999999999
$ws.Range("L134").Value = 12461.5005
$ws.Range("M134").Value = -5229.545999999999
$ws.Range("N134").Value = -17531.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 526677.8
$ws.Range("J92").Value = 408.72726
$ws.Range("L92").Value = 1226.18178
$ws.Range("N92").Value = -3722.18178
$ws.Range("H131").Value = 729.76117
$ws.Range("J131").Value = 911.02325
$ws.Range("L131").Value = 2733.06975
$ws.Range("N131").Value = -12813.06975
$ws.Range("H132").Value = 843271.7
$ws.Range("I132").Value = 1192.5
$ws.Range("J132").Value = 1264311.2
$ws.Range("K132").Value = 10732.5
$ws.Range("L132").Value = 11378800.8
$ws.Range("M132").Value = -8202.5
$ws.Range("N132").Value = -11383860.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1246.75
$ws.Range("I97").Value = 1275
$ws.Range("J97").Value = 1105.5
$ws.Range("K97").Value = 1275
$ws.Range("L97").Value = 1105.5
$ws.Range("M97").Value = -779
$ws.Range("N97").Value = -2097.5
$ws.Range("H102").Value = 1872.32
$ws.Range("I102").Value = 1879
$ws.Range("J102").Value = 1823.3334
$ws.Range("K102").Value = 1879
$ws.Range("L102").Value = 1823.3334
$ws.Range("M102").Value = -257
$ws.Range("N102").Value = -5067.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5164.6523
$ws.Range("I93").Value = 6691.5293
$ws.Range("J93").Value = 838.5
$ws.Range("K93").Value = 6691.5293
$ws.Range("L93").Value = 838.5
$ws.Range("M93").Value = -5443.5293
$ws.Range("N93").Value = -3334.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9093409
$ws.Range("I62").Value = 25001650
$ws.Range("J62").Value = 2985.7144
$ws.Range("K62").Value = 25001650
$ws.Range("L62").Value = 2985.7144
$ws.Range("M62").Value = -25001026
$ws.Range("N62").Value = -4233.7144
$ws.Range("H65").Value = 9093409
$ws.Range("I65").Value = 25001650
$ws.Range("J65").Value = 2985.7144
$ws.Range("K65").Value = 125008250
$ws.Range("L65").Value = 14928.572
$ws.Range("M65").Value = -125005130
$ws.Range("N65").Value = -21168.572
$ws.Range("H96").Value = 2560.25
$ws.Range("I96").Value = 1971.2858
$ws.Range("J96").Value = 3384.8
$ws.Range("K96").Value = 1971.2858
$ws.Range("L96").Value = 3384.8
$ws.Range("M96").Value = -598.2858000000001
$ws.Range("N96").Value = -6130.8
$ws.Range("H107").Value = 20834182
$ws.Range("I107").Value = 33333914
$ws.Range("J107").Value = 1293
$ws.Range("K107").Value = 100001742
$ws.Range("L107").Value = 3879
$ws.Range("M107").Value = -99999822
$ws.Range("N107").Value = -7719
$ws.Range("H113").Value = 606
$ws.Range("I113").Value = 319.66666
$ws.Range("K113").Value = 958.9999799999999
$ws.Range("M113").Value = 1211.00002
$ws.Range("H131").Value = 21863
$ws.Range("J131").Value = 21863
$ws.Range("L131").Value = 21863
$ws.Range("N131").Value = -31943
